# Fill in row 12 (2025-10-22) data on Sheet1 and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12").Value = 45952
$ws.Range("B12").Value = 5598
$ws.Range("C12").Value = 4374
$ws.Range("D12").Value = 4062
$ws.Range("E12").Value = 239
$ws.Range("F12").Value = 38
$ws.Range("G12").Value = 31
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 1

$ws.Range("F21").Select()
